# Update currentAveragePrice / profit columns with refreshed market data
# (values captured by the scheduled market-data runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 465.3
$ws.Cells.Item(29, 9).Value = 225.75
$ws.Cells.Item(29, 10).Value = 625
$ws.Cells.Item(29, 11).Value = 677.25
$ws.Cells.Item(29, 12).Value = 1875
$ws.Cells.Item(29, 13).Value = -396.25
$ws.Cells.Item(29, 14).Value = -2437

$ws.Cells.Item(38, 8).Value = 2501.4546
$ws.Cells.Item(38, 9).Value = 130
$ws.Cells.Item(38, 10).Value = 4143.231
$ws.Cells.Item(38, 11).Value = 390
$ws.Cells.Item(38, 12).Value = 12429.693
$ws.Cells.Item(38, 13).Value = -18
$ws.Cells.Item(38, 14).Value = -13173.693

$ws.Cells.Item(58, 8).Value = 1714.2916
$ws.Cells.Item(58, 9).Value = 53.75
$ws.Cells.Item(58, 10).Value = 2544.5625
$ws.Cells.Item(58, 11).Value = 161.25
$ws.Cells.Item(58, 12).Value = 7633.6875
$ws.Cells.Item(58, 13).Value = -11.25
$ws.Cells.Item(58, 14).Value = -7933.6875

$ws.Cells.Item(87, 8).Value = 31836.115
$ws.Cells.Item(87, 10).Value = 31836.115
$ws.Cells.Item(87, 12).Value = 31836.115
$ws.Cells.Item(87, 14).Value = -34332.11500000001

$ws.Cells.Item(90, 8).Value = 31836.115
$ws.Cells.Item(90, 10).Value = 31836.115
$ws.Cells.Item(90, 12).Value = 95508.345
$ws.Cells.Item(90, 14).Value = -107988.345

$ws.Cells.Item(129, 8).Value = 1096.2456
$ws.Cells.Item(129, 10).Value = 1131.2222
$ws.Cells.Item(129, 12).Value = 3393.6666
$ws.Cells.Item(129, 14).Value = -13393.6666

$ws.Cells.Item(135, 8).Value = 26848.564
$ws.Cells.Item(135, 9).Value = 28683.166
$ws.Cells.Item(135, 10).Value = 4833.3335
$ws.Cells.Item(135, 11).Value = 258148.494
$ws.Cells.Item(135, 12).Value = 43500.0015
$ws.Cells.Item(135, 13).Value = -255613.494
$ws.Cells.Item(135, 14).Value = -48570.0015

$ws.Cells.Item(138, 8).Value = 373087.88
$ws.Cells.Item(138, 9).Value = 2164.65
$ws.Cells.Item(138, 10).Value = 836741.9399999999
$ws.Cells.Item(138, 11).Value = 6493.950000000001
$ws.Cells.Item(138, 12).Value = 2510225.82
$ws.Cells.Item(138, 13).Value = -1353.950000000001
$ws.Cells.Item(138, 14).Value = -2520505.82


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 25147.777
$ws.Cells.Item(32, 9).Value = 4504.346
$ws.Cells.Item(32, 10).Value = 122734.91
$ws.Cells.Item(32, 11).Value = 4504.346
$ws.Cells.Item(32, 12).Value = 122734.91
$ws.Cells.Item(32, 13).Value = -4217.346
$ws.Cells.Item(32, 14).Value = -123308.91

$ws.Cells.Item(63, 8).Value = 2037.4286
$ws.Cells.Item(63, 9).Value = 1897.6364
$ws.Cells.Item(63, 10).Value = 2550
$ws.Cells.Item(63, 11).Value = 1897.6364
$ws.Cells.Item(63, 12).Value = 2550
$ws.Cells.Item(63, 13).Value = -1211.6364
$ws.Cells.Item(63, 14).Value = -3922

$ws.Cells.Item(66, 8).Value = 2037.4286
$ws.Cells.Item(66, 9).Value = 1897.6364
$ws.Cells.Item(66, 10).Value = 2550
$ws.Cells.Item(66, 11).Value = 9488.182000000001
$ws.Cells.Item(66, 12).Value = 12750
$ws.Cells.Item(66, 13).Value = -6056.182000000001
$ws.Cells.Item(66, 14).Value = -19614


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(33, 8).Value = 33640.332
$ws.Cells.Item(33, 9).Value = 33640.332
$ws.Cells.Item(33, 11).Value = 33640.332
$ws.Cells.Item(33, 13).Value = -33304.332


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 6787.5
$ws.Cells.Item(35, 9).Value = 1850
$ws.Cells.Item(35, 11).Value = 1850
$ws.Cells.Item(35, 13).Value = -1556

$ws.Cells.Item(86, 8).Value = 90919180
$ws.Cells.Item(86, 9).Value = 166682190
$ws.Cells.Item(86, 10).Value = 3580
$ws.Cells.Item(86, 11).Value = 166682190
$ws.Cells.Item(86, 12).Value = 3580
$ws.Cells.Item(86, 13).Value = -166681067
$ws.Cells.Item(86, 14).Value = -5826

$ws.Cells.Item(89, 8).Value = 90919180
$ws.Cells.Item(89, 9).Value = 166682190
$ws.Cells.Item(89, 10).Value = 3580
$ws.Cells.Item(89, 11).Value = 833410950
$ws.Cells.Item(89, 12).Value = 17900
$ws.Cells.Item(89, 13).Value = -833405334
$ws.Cells.Item(89, 14).Value = -29132

$ws.Cells.Item(99, 8).Value = 1605.5555
$ws.Cells.Item(99, 9).Value = 1421.4286
$ws.Cells.Item(99, 11).Value = 1421.4286
$ws.Cells.Item(99, 13).Value = 76.57140000000004

$ws.Cells.Item(126, 8).Value = 1605.5555
$ws.Cells.Item(126, 9).Value = 1421.4286
$ws.Cells.Item(126, 11).Value = 4264.2858
$ws.Cells.Item(126, 13).Value = -1794.2858


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 329.89655
$ws.Cells.Item(5, 9).Value = 316.67856
$ws.Cells.Item(5, 10).Value = 700
$ws.Cells.Item(5, 11).Value = 950.03568
$ws.Cells.Item(5, 12).Value = 2100
$ws.Cells.Item(5, 13).Value = -838.03568
$ws.Cells.Item(5, 14).Value = -2324

$ws.Cells.Item(34, 8).Value = 1963.0769
$ws.Cells.Item(34, 9).Value = 544
$ws.Cells.Item(34, 10).Value = 2850
$ws.Cells.Item(34, 11).Value = 1632
$ws.Cells.Item(34, 12).Value = 8550
$ws.Cells.Item(34, 13).Value = -1548
$ws.Cells.Item(34, 14).Value = -8718

$ws.Cells.Item(39, 8).Value = 1604.5667
$ws.Cells.Item(39, 10).Value = 1604.5667
$ws.Cells.Item(39, 12).Value = 4813.7001
$ws.Cells.Item(39, 14).Value = -5401.7001

$ws.Cells.Item(55, 8).Value = 42510.56
$ws.Cells.Item(55, 9).Value = 733.3333
$ws.Cells.Item(55, 10).Value = 48207.453
$ws.Cells.Item(55, 11).Value = 2199.9999
$ws.Cells.Item(55, 12).Value = 144622.359
$ws.Cells.Item(55, 13).Value = -2022.9999
$ws.Cells.Item(55, 14).Value = -144976.359

$ws.Cells.Item(107, 8).Value = 1177.7142
$ws.Cells.Item(107, 9).Value = 855.94446
$ws.Cells.Item(107, 10).Value = 1419.0416
$ws.Cells.Item(107, 11).Value = 2567.83338
$ws.Cells.Item(107, 12).Value = 4257.1248
$ws.Cells.Item(107, 13).Value = -647.83338
$ws.Cells.Item(107, 14).Value = -8097.1248

$ws.Cells.Item(113, 8).Value = 514.0454999999999
$ws.Cells.Item(113, 9).Value = 470.08334
$ws.Cells.Item(113, 10).Value = 566.8
$ws.Cells.Item(113, 11).Value = 1410.25002
$ws.Cells.Item(113, 12).Value = 1700.4
$ws.Cells.Item(113, 13).Value = 759.7499800000001
$ws.Cells.Item(113, 14).Value = -6040.4

$ws.Cells.Item(129, 8).Value = 54053.316
$ws.Cells.Item(129, 10).Value = 1981.4445
$ws.Cells.Item(129, 12).Value = 5944.333500000001
$ws.Cells.Item(129, 14).Value = -15944.3335

$ws.Cells.Item(131, 8).Value = 812.1818
$ws.Cells.Item(131, 9).Value = 403.16666
$ws.Cells.Item(131, 10).Value = 1303
$ws.Cells.Item(131, 11).Value = 1209.49998
$ws.Cells.Item(131, 12).Value = 3909
$ws.Cells.Item(131, 13).Value = 3830.50002
$ws.Cells.Item(131, 14).Value = -13989

$ws.Cells.Item(132, 8).Value = 881418.25
$ws.Cells.Item(132, 9).Value = 1647160.6
$ws.Cells.Item(132, 10).Value = 6284.143
$ws.Cells.Item(132, 11).Value = 14824445.4
$ws.Cells.Item(132, 12).Value = 56557.287
$ws.Cells.Item(132, 13).Value = -14821915.4
$ws.Cells.Item(132, 14).Value = -61617.287

$ws.Cells.Item(135, 8).Value = 329.89655
$ws.Cells.Item(135, 9).Value = 316.67856
$ws.Cells.Item(135, 10).Value = 700
$ws.Cells.Item(135, 11).Value = 2850.10704
$ws.Cells.Item(135, 12).Value = 6300
$ws.Cells.Item(135, 13).Value = -315.1070399999999
$ws.Cells.Item(135, 14).Value = -11370


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 700
$ws.Cells.Item(29, 9).Value = 700
$ws.Cells.Item(29, 11).Value = 700
$ws.Cells.Item(29, 13).Value = -410

$ws.Cells.Item(102, 8).Value = 1904.7222
$ws.Cells.Item(102, 9).Value = 1560.3846
$ws.Cells.Item(102, 10).Value = 2800
$ws.Cells.Item(102, 11).Value = 1560.3846
$ws.Cells.Item(102, 12).Value = 2800
$ws.Cells.Item(102, 13).Value = 61.61539999999991
$ws.Cells.Item(102, 14).Value = -6044

$ws.Cells.Item(107, 8).Value = 887.5599999999999
$ws.Cells.Item(107, 9).Value = 791
$ws.Cells.Item(107, 10).Value = 1010.4545
$ws.Cells.Item(107, 11).Value = 791
$ws.Cells.Item(107, 12).Value = 1010.4545
$ws.Cells.Item(107, 13).Value = 1129
$ws.Cells.Item(107, 14).Value = -4850.4545

$ws.Cells.Item(132, 8).Value = 2768.913
$ws.Cells.Item(132, 9).Value = 2437.2856
$ws.Cells.Item(132, 10).Value = 3284.7778
$ws.Cells.Item(132, 11).Value = 7311.8568
$ws.Cells.Item(132, 12).Value = 9854.3334
$ws.Cells.Item(132, 13).Value = -4781.8568
$ws.Cells.Item(132, 14).Value = -14914.3334


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1134.6154
$ws.Cells.Item(22, 9).Value = 950
$ws.Cells.Item(22, 10).Value = 1216.6666
$ws.Cells.Item(22, 11).Value = 950
$ws.Cells.Item(22, 12).Value = 1216.6666
$ws.Cells.Item(22, 13).Value = -655
$ws.Cells.Item(22, 14).Value = -1806.6666

$ws.Cells.Item(27, 8).Value = 1134.6154
$ws.Cells.Item(27, 9).Value = 950
$ws.Cells.Item(27, 10).Value = 1216.6666
$ws.Cells.Item(27, 11).Value = 950
$ws.Cells.Item(27, 12).Value = 1216.6666
$ws.Cells.Item(27, 13).Value = -843
$ws.Cells.Item(27, 14).Value = -1430.6666

$ws.Cells.Item(46, 8).Value = 1501.1818
$ws.Cells.Item(46, 9).Value = 1222.5555
$ws.Cells.Item(46, 10).Value = 2755
$ws.Cells.Item(46, 11).Value = 1222.5555
$ws.Cells.Item(46, 12).Value = 2755
$ws.Cells.Item(46, 13).Value = -1034.5555
$ws.Cells.Item(46, 14).Value = -3131


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 7987.125
$ws.Cells.Item(81, 9).Value = 17985.334
$ws.Cells.Item(81, 10).Value = 1988.2
$ws.Cells.Item(81, 11).Value = 35970.668
$ws.Cells.Item(81, 12).Value = 3976.4
$ws.Cells.Item(81, 13).Value = -34909.668
$ws.Cells.Item(81, 14).Value = -6098.4

$ws.Cells.Item(84, 8).Value = 7987.125
$ws.Cells.Item(84, 9).Value = 17985.334
$ws.Cells.Item(84, 10).Value = 1988.2
$ws.Cells.Item(84, 11).Value = 179853.34
$ws.Cells.Item(84, 12).Value = 19882
$ws.Cells.Item(84, 13).Value = -174549.34
$ws.Cells.Item(84, 14).Value = -30490

$ws.Cells.Item(96, 8).Value = 1866.6666
$ws.Cells.Item(96, 9).Value = 1950
$ws.Cells.Item(96, 10).Value = 1700
$ws.Cells.Item(96, 11).Value = 1950
$ws.Cells.Item(96, 12).Value = 1700
$ws.Cells.Item(96, 13).Value = -577
$ws.Cells.Item(96, 14).Value = -4446

